# cv124013a.xlsx - correção nos dados e inicio da analise PNAD 2009
#
# 1. Rename the mislabeled column header "unnamed: 1_level_1" -> "total"
# 2. Remove the two sub-heading-only rows that never carried any data
#    ("situação do domicílio" and "grandes regiões e unidades da federação").
#    Deleting them shifts every following row (and its data) up by two,
#    so the table now ends two rows earlier, at row 38 (distrito federal)
#    instead of row 40 - matching the new dimension A1:H38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the column header text.
$ws.Range("B2").Value = "total"

# Delete the "situação do domicílio" sub-header row (row 5).
$ws.Rows("5:5").Delete()

# After the row above was removed, the "grandes regiões e unidades da
# federação" sub-header row (originally row 8) is now row 7.
$ws.Rows("7:7").Delete()
